$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.868554723999523
$ws.Cells.Item(2, 3).Value = 0.2230895546694285
$ws.Cells.Item(2, 4).Value = 0.02541487974610135
$ws.Cells.Item(2, 5).Value = 0.1093082662282425
$ws.Cells.Item(2, 6).Value = 0.7337237621415937
$ws.Cells.Item(2, 8).Value = 0.07973214163530429
$ws.Cells.Item(2, 9).Value = 0.6503881080935159
$ws.Cells.Item(2, 12).Value = 0.2035184872959732
$ws.Cells.Item(2, 13).Value = 0.2039161632910051
$ws.Cells.Item(2, 14).Value = 1.278446630484339
$ws.Cells.Item(2, 15).Value = 2.510429043626402

$ws.Cells.Item(3, 2).Value = 0.7866008231924866
$ws.Cells.Item(3, 3).Value = 0.209625235377473
$ws.Cells.Item(3, 4).Value = 0.02394862812381859
$ws.Cells.Item(3, 5).Value = 0.1102144949870643
$ws.Cells.Item(3, 6).Value = 0.7287116675615835
$ws.Cells.Item(3, 8).Value = 0.07973214163530429
$ws.Cells.Item(3, 9).Value = 0.6558186271010094
$ws.Cells.Item(3, 12).Value = 0.200750244957085
$ws.Cells.Item(3, 13).Value = 0.1906489805552667
$ws.Cells.Item(3, 14).Value = 1.289411023549796
$ws.Cells.Item(3, 15).Value = 2.507811326002326

$ws.Cells.Item(4, 2).Value = 0.7363698976328124
$ws.Cells.Item(4, 3).Value = 0.2012869645223674
$ws.Cells.Item(4, 4).Value = 0.02304062198928847
$ws.Cells.Item(4, 5).Value = 0.1108039141966772
$ws.Cells.Item(4, 6).Value = 0.7261286202203507
$ws.Cells.Item(4, 8).Value = 0.07973214163530429
$ws.Cells.Item(4, 9).Value = 0.6595363778697489
$ws.Cells.Item(4, 12).Value = 0.1991511633465564
$ws.Cells.Item(4, 13).Value = 0.1825620133722836
$ws.Cells.Item(4, 14).Value = 1.296623658830342
$ws.Cells.Item(4, 15).Value = 2.507846823385421

$ws.Cells.Item(5, 2).Value = 0.7159241101786051
$ws.Cells.Item(5, 3).Value = 0.197871324049089
$ws.Cells.Item(5, 4).Value = 0.02266868007407652
$ws.Cells.Item(5, 5).Value = 0.1110524143337099
$ws.Cells.Item(5, 6).Value = 0.7252002856644708
$ws.Cells.Item(5, 8).Value = 0.07973214163530429
$ws.Cells.Item(5, 9).Value = 0.6611477443888951
$ws.Cells.Item(5, 12).Value = 0.1985248895987368
$ws.Cells.Item(5, 13).Value = 0.1792815925659283
$ws.Cells.Item(5, 14).Value = 1.299683767271169
$ws.Cells.Item(5, 15).Value = 2.508274158580491

$ws.Cells.Item(6, 2).Value = 0.7125305743480794
$ws.Cells.Item(6, 3).Value = 0.1973030934977089
$ws.Cells.Item(6, 4).Value = 0.0226068038938152
$ws.Cells.Item(6, 5).Value = 0.1110941796437102
$ws.Cells.Item(6, 6).Value = 0.7250536412817326
$ws.Cells.Item(6, 8).Value = 0.07973214163530429
$ws.Cells.Item(6, 9).Value = 0.6614211284807361
$ws.Cells.Item(6, 12).Value = 0.1984224310905489
$ws.Cells.Item(6, 13).Value = 0.1787377985183909
$ws.Cells.Item(6, 14).Value = 1.300199200698486
$ws.Cells.Item(6, 15).Value = 2.508370047095468

$ws.Cells.Item(7, 2).Value = 0.7360940607944144
$ws.Cells.Item(7, 3).Value = 0.2012409715687511
$ws.Cells.Item(7, 4).Value = 0.02303561360680817
$ws.Cells.Item(7, 5).Value = 0.1108072319065438
$ws.Cells.Item(7, 6).Value = 0.7261155972353777
$ws.Cells.Item(7, 8).Value = 0.07973214163530429
$ws.Cells.Item(7, 9).Value = 0.6595577192648285
$ws.Cells.Item(7, 12).Value = 0.1991426144169424
$ws.Cells.Item(7, 13).Value = 0.1825177111122684
$ws.Cells.Item(7, 14).Value = 1.296664438882367
$ws.Cells.Item(7, 15).Value = 2.507850915182161

$ws.Cells.Item(8, 2).Value = 0.8402792903086151
$ws.Cells.Item(8, 3).Value = 0.2184619098921132
$ws.Cells.Item(8, 4).Value = 0.02491093020738333
$ws.Cells.Item(8, 5).Value = 0.1096138933733757
$ws.Cells.Item(8, 6).Value = 0.7318929501761104
$ws.Cells.Item(8, 8).Value = 0.07973214163530429
$ws.Cells.Item(8, 9).Value = 0.6521809458837176
$ws.Cells.Item(8, 12).Value = 0.2025431465326761
$ws.Cells.Item(8, 13).Value = 0.1993294885757422
$ws.Cells.Item(8, 14).Value = 1.282127486823981
$ws.Cells.Item(8, 15).Value = 2.509185352070745

$ws.Cells.Item(9, 2).Value = 1.04524546286973
$ws.Cells.Item(9, 3).Value = 0.2516623280749002
$ws.Cells.Item(9, 4).Value = 0.02852644670926452
$ws.Cells.Item(9, 5).Value = 0.1075350074939618
$ws.Cells.Item(9, 6).Value = 0.7471488321144761
$ws.Cells.Item(9, 8).Value = 0.07973214163530429
$ws.Cells.Item(9, 9).Value = 0.6407597876439617
$ws.Cells.Item(9, 12).Value = 0.2100079451870585
$ws.Cells.Item(9, 13).Value = 0.2327589209310261
$ws.Cells.Item(9, 14).Value = 1.257428932265888
$ws.Cells.Item(9, 15).Value = 2.524851093627774

$ws.Cells.Item(10, 2).Value = 1.196186840197811
$ws.Cells.Item(10, 3).Value = 0.2757021534029889
$ws.Cells.Item(10, 4).Value = 0.03114427024885202
$ws.Cells.Item(10, 5).Value = 0.1061661635196298
$ws.Cells.Item(10, 6).Value = 0.7607586634862429
$ws.Cells.Item(10, 8).Value = 0.07973214163530429
$ws.Cells.Item(10, 9).Value = 0.6342287566908382
$ws.Cells.Item(10, 12).Value = 0.2159758550657216
$ws.Cells.Item(10, 13).Value = 0.2575932380129018
$ws.Cells.Item(10, 14).Value = 1.241599453196393
$ws.Cells.Item(10, 15).Value = 2.544341185788653

$ws.Cells.Item(11, 2).Value = 1.264921059430606
$ws.Cells.Item(11, 3).Value = 0.2865609784079766
$ws.Cells.Item(11, 4).Value = 0.03232668827989471
$ws.Cells.Item(11, 5).Value = 0.1055777004257199
$ws.Cells.Item(11, 6).Value = 0.7674733401018159
$ws.Cells.Item(11, 8).Value = 0.07973214163530429
$ws.Cells.Item(11, 9).Value = 0.6316624086232991
$ws.Cells.Item(11, 12).Value = 0.2187954406617934
$ws.Cells.Item(11, 13).Value = 0.2689489737203914
$ws.Cells.Item(11, 14).Value = 1.234900224204125
$ws.Cells.Item(11, 15).Value = 2.554946461127571

$ws.Cells.Item(12, 2).Value = 1.290957879753762
$ws.Cells.Item(12, 3).Value = 0.2906617261323561
$ws.Cells.Item(12, 4).Value = 0.03277320792186345
$ws.Cells.Item(12, 5).Value = 0.105359775680931
$ws.Cells.Item(12, 6).Value = 0.7700913820657007
$ws.Cells.Item(12, 8).Value = 0.07973214163530429
$ws.Cells.Item(12, 9).Value = 0.6307488579327938
$ws.Cells.Item(12, 12).Value = 0.2198781590972061
$ws.Cells.Item(12, 13).Value = 0.2732573227625466
$ws.Cells.Item(12, 14).Value = 1.232435478975681
$ws.Cells.Item(12, 15).Value = 2.559212844790693

$ws.Cells.Item(13, 2).Value = 1.285350022181035
$ws.Cells.Item(13, 3).Value = 0.289779060101381
$ws.Cells.Item(13, 4).Value = 0.03267709727199986
$ws.Cells.Item(13, 5).Value = 0.10540649128765
$ws.Cells.Item(13, 6).Value = 0.7695241881704931
$ws.Cells.Item(13, 8).Value = 0.07973214163530429
$ws.Cells.Item(13, 9).Value = 0.6309430145696737
$ws.Cells.Item(13, 12).Value = 0.2196443100375234
$ws.Cells.Item(13, 13).Value = 0.2723290826513036
$ws.Cells.Item(13, 14).Value = 1.232963099458267
$ws.Cells.Item(13, 15).Value = 2.558282861389472

$ws.Cells.Item(14, 2).Value = 1.267062958902045
$ws.Cells.Item(14, 3).Value = 0.2868985759561724
$ws.Cells.Item(14, 4).Value = 0.03236344862632023
$ws.Cells.Item(14, 5).Value = 0.1055596732035387
$ws.Cells.Item(14, 6).Value = 0.7676872177237897
$ws.Cells.Item(14, 8).Value = 0.07973214163530429
$ws.Cells.Item(14, 9).Value = 0.631586081983059
$ws.Cells.Item(14, 12).Value = 0.2188842162342439
$ws.Cells.Item(14, 13).Value = 0.269303261604847
$ws.Cells.Item(14, 14).Value = 1.234696003169653
$ws.Cells.Item(14, 15).Value = 2.555292439001732

$ws.Cells.Item(15, 2).Value = 1.25586269024285
$ws.Cells.Item(15, 3).Value = 0.2851327253684417
$ws.Cells.Item(15, 4).Value = 0.03217116809584297
$ws.Cells.Item(15, 5).Value = 0.105654141164234
$ws.Cells.Item(15, 6).Value = 0.7665718335318132
$ws.Cells.Item(15, 8).Value = 0.07973214163530429
$ws.Cells.Item(15, 9).Value = 0.6319875701771025
$ws.Cells.Item(15, 12).Value = 0.2184205888016493
$ws.Cells.Item(15, 13).Value = 0.2674509170577579
$ws.Cells.Item(15, 14).Value = 1.235766846310739
$ws.Cells.Item(15, 15).Value = 2.55349333706323

$ws.Cells.Item(16, 2).Value = 1.19169618775328
$ws.Cells.Item(16, 3).Value = 0.2749909386322997
$ws.Cells.Item(16, 4).Value = 0.03106682464807875
$ws.Cells.Item(16, 5).Value = 0.1062053090172772
$ws.Cells.Item(16, 6).Value = 0.7603303829543648
$ws.Cells.Item(16, 8).Value = 0.07973214163530429
$ws.Cells.Item(16, 9).Value = 0.6344046231386997
$ws.Cells.Item(16, 12).Value = 0.2157936912231833
$ws.Cells.Item(16, 13).Value = 0.2568522710006107
$ws.Cells.Item(16, 14).Value = 1.242047355693195
$ws.Cells.Item(16, 15).Value = 2.543683128564908

$ws.Cells.Item(17, 2).Value = 1.152349082233172
$ws.Cells.Item(17, 3).Value = 0.2687494286349761
$ws.Cells.Item(17, 4).Value = 0.03038716738219449
$ws.Cells.Item(17, 5).Value = 0.1065521936406091
$ws.Cells.Item(17, 6).Value = 0.7566355793734658
$ws.Cells.Item(17, 8).Value = 0.07973214163530429
$ws.Cells.Item(17, 9).Value = 0.6359911055545631
$ws.Cells.Item(17, 12).Value = 0.2142089612493265
$ws.Cells.Item(17, 13).Value = 0.2503651587298208
$ws.Cells.Item(17, 14).Value = 1.246028711977139
$ws.Cells.Item(17, 15).Value = 2.538110539547318

$ws.Cells.Item(18, 2).Value = 1.129724363880428
$ws.Cells.Item(18, 3).Value = 0.2651522391398373
$ws.Cells.Item(18, 4).Value = 0.02999545313528529
$ws.Cells.Item(18, 5).Value = 0.1067549350036736
$ws.Cells.Item(18, 6).Value = 0.7545596964181698
$ws.Cells.Item(18, 8).Value = 0.07973214163530429
$ws.Cells.Item(18, 9).Value = 0.636941691065541
$ws.Cells.Item(18, 12).Value = 0.2133073300191057
$ws.Cells.Item(18, 13).Value = 0.2466394615731815
$ws.Cells.Item(18, 14).Value = 1.248365905097849
$ws.Cells.Item(18, 15).Value = 2.535069010940532

$ws.Cells.Item(19, 2).Value = 1.122065216821397
$ws.Cells.Item(19, 3).Value = 0.2639330549127124
$ws.Cells.Item(19, 4).Value = 0.02986268984005847
$ws.Cells.Item(19, 5).Value = 0.1068241334857625
$ws.Cells.Item(19, 6).Value = 0.7538652984473231
$ws.Cells.Item(19, 8).Value = 0.07973214163530429
$ws.Cells.Item(19, 9).Value = 0.6372700812599632
$ws.Cells.Item(19, 12).Value = 0.2130037491011478
$ws.Cells.Item(19, 13).Value = 0.2453789593374296
$ws.Cells.Item(19, 14).Value = 1.249165350095794
$ws.Cells.Item(19, 15).Value = 2.534067303822383

$ws.Cells.Item(20, 2).Value = 1.156536966175679
$ws.Cells.Item(20, 3).Value = 0.2694145984389138
$ws.Cells.Item(20, 4).Value = 0.03045960032996931
$ws.Cells.Item(20, 5).Value = 0.1065149337281874
$ws.Cells.Item(20, 6).Value = 0.7570237977754317
$ws.Cells.Item(20, 8).Value = 0.07973214163530429
$ws.Cells.Item(20, 9).Value = 0.6358182797765402
$ws.Cells.Item(20, 12).Value = 0.2143766380621912
$ws.Cells.Item(20, 13).Value = 0.2510551534617562
$ws.Cells.Item(20, 14).Value = 1.245600002728928
$ws.Cells.Item(20, 15).Value = 2.538686808921028

$ws.Cells.Item(21, 2).Value = 1.272434088987893
$ws.Cells.Item(21, 3).Value = 0.287744950823452
$ws.Cells.Item(21, 4).Value = 0.03245560857011043
$ws.Cells.Item(21, 5).Value = 0.1055145466825134
$ws.Cells.Item(21, 6).Value = 0.7682247352931881
$ws.Cells.Item(21, 8).Value = 0.07973214163530429
$ws.Cells.Item(21, 9).Value = 0.6313956154277349
$ws.Cells.Item(21, 12).Value = 0.2191070675293787
$ws.Cells.Item(21, 13).Value = 0.2701917984746984
$ws.Cells.Item(21, 14).Value = 1.234185050618017
$ws.Cells.Item(21, 15).Value = 2.556164000409723

$ws.Cells.Item(22, 2).Value = 1.348229389395613
$ws.Cells.Item(22, 3).Value = 0.2996592177675836
$ws.Cells.Item(22, 4).Value = 0.03375289826872319
$ws.Cells.Item(22, 5).Value = 0.1048893712228889
$ws.Cells.Item(22, 6).Value = 0.7759843601248235
$ws.Cells.Item(22, 8).Value = 0.07973214163530429
$ws.Cells.Item(22, 9).Value = 0.6288448261447215
$ws.Cells.Item(22, 12).Value = 0.2222860947501886
$ws.Cells.Item(22, 13).Value = 0.2827462649922197
$ws.Cells.Item(22, 14).Value = 1.22714500381219
$ws.Cells.Item(22, 15).Value = 2.569045997061778

$ws.Cells.Item(23, 2).Value = 1.307771925503005
$ws.Cells.Item(23, 3).Value = 0.2933064180554368
$ws.Cells.Item(23, 4).Value = 0.03306117807824194
$ws.Cells.Item(23, 5).Value = 0.1052204219239066
$ws.Cells.Item(23, 6).Value = 0.7718026987780462
$ws.Cells.Item(23, 8).Value = 0.07973214163530429
$ws.Cells.Item(23, 9).Value = 0.6301751230646033
$ws.Cells.Item(23, 12).Value = 0.2205814094638896
$ws.Cells.Item(23, 13).Value = 0.2760414342386284
$ws.Cells.Item(23, 14).Value = 1.230863963343076
$ws.Cells.Item(23, 15).Value = 2.562036973912711

$ws.Cells.Item(24, 2).Value = 1.154643634772867
$ws.Cells.Item(24, 3).Value = 0.2691139027641043
$ws.Cells.Item(24, 4).Value = 0.03042685641456444
$ws.Cells.Item(24, 5).Value = 0.1065317686151008
$ws.Cells.Item(24, 6).Value = 0.7568481337617641
$ws.Cells.Item(24, 8).Value = 0.07973214163530429
$ws.Cells.Item(24, 9).Value = 0.6358962944046027
$ws.Cells.Item(24, 12).Value = 0.2143008019317563
$ws.Cells.Item(24, 13).Value = 0.2507431948956977
$ws.Cells.Item(24, 14).Value = 1.245793671853825
$ws.Cells.Item(24, 15).Value = 2.538425772231165

$ws.Cells.Item(25, 2).Value = 0.9897311952612995
$ws.Cells.Item(25, 3).Value = 0.2427422693882022
$ws.Cells.Item(25, 4).Value = 0.02755506235195782
$ws.Cells.Item(25, 5).Value = 0.1080695106552074
$ws.Cells.Item(25, 6).Value = 0.7426005586038329
$ws.Cells.Item(25, 8).Value = 0.07973214163530429
$ws.Cells.Item(25, 9).Value = 0.6435231486489457
$ws.Cells.Item(25, 12).Value = 0.207903465731647
$ws.Cells.Item(25, 13).Value = 0.2236667525888052
$ws.Cells.Item(25, 14).Value = 1.263703380562752
$ws.Cells.Item(25, 15).Value = 2.519213474105328
